# Update the portfolio workbook:
#  - Sheet "交易记录" (transactions): re-date existing rows to 2025-10-14 and
#    append 6 new "买入" (buy) rows (12-17).
#  - Sheet "当前持仓" (current holdings): refresh cash / prices / P&L figures
#    and the two date columns, and append 4 new holding rows (7-10).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 交易记录 (transaction log)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Dates are stored as plain text (not real Excel dates), so force the
# "Text" number format before writing date-looking strings — otherwise
# Excel's COM layer auto-coerces "2025-10-14" into a date serial.
$ws1.Range("A2:A17").NumberFormat = "@"

# --- Re-date the existing rows (values otherwise unchanged) ---
$ws1.Cells.Item(2, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(3, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(4, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(5, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(6, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(7, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(8, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(9, 1).Value()  = "2025-10-14"
$ws1.Cells.Item(10, 1).Value() = "2025-10-14"
$ws1.Cells.Item(11, 1).Value() = "2025-10-14"

# --- New transaction rows ---
function Set-TxRow($row, $date, $op, $asset, $qty, $price, $amount, $fee, $cost) {
    $ws1.Cells.Item($row, 1).Value() = $date
    $ws1.Cells.Item($row, 2).Value() = $op
    $ws1.Cells.Item($row, 3).Value() = $asset
    $ws1.Cells.Item($row, 4).Value() = $qty
    $ws1.Cells.Item($row, 5).Value() = $price
    $ws1.Cells.Item($row, 6).Value() = $amount
    $ws1.Cells.Item($row, 7).Value() = $fee
    $ws1.Cells.Item($row, 8).Value() = $cost
}

Set-TxRow 12 "2025-10-14" "买入" "有色etf" 1800 0.901 1621.8 5 1626.8
Set-TxRow 13 "2025-10-14" "买入" "芯片etf" 4000 1.1   4400   5 4405
Set-TxRow 14 "2025-10-14" "买入" "科创50"   1500 1.532 2298   5 2303
Set-TxRow 15 "2025-10-14" "买入" "芯片etf" 1500 1.055 1582.5 5 1587.5
Set-TxRow 16 "2025-10-14" "买入" "科创50"   1000 1.497 1497   5 1502
Set-TxRow 17 "2025-10-14" "买入" "稀土etf" 800  1.355 1084   5 1089

# ---------------------------------------------------------------------
# Sheet 2: 当前持仓 (current holdings)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Same text-date handling as sheet 1, for columns I (price update date)
# and J (last update date), rows 2-10.
$ws2.Range("I2:J10").NumberFormat = "@"

# --- Row 2: 现金 (cash) ---
$ws2.Cells.Item(2, 2).Value()  = 18711.377216
$ws2.Cells.Item(2, 3).Value()  = 1
$ws2.Cells.Item(2, 4).Value()  = 18711.377216
$ws2.Cells.Item(2, 5).Value()  = 1
$ws2.Cells.Item(2, 6).Value()  = 18711.377216
$ws2.Cells.Item(2, 7).Value()  = 0
$ws2.Cells.Item(2, 8).Value()  = 0
$ws2.Cells.Item(2, 9).Value()  = "2025-10-14"
$ws2.Cells.Item(2, 10).Value() = "2025-10-14"

# --- Row 3: 华泰柏瑞中证红利低波动ETF链接A ---
$ws2.Cells.Item(3, 5).Value()  = 1.6741
$ws2.Cells.Item(3, 6).Value()  = 20725.039921
$ws2.Cells.Item(3, 7).Value()  = -283.6322121521698
$ws2.Cells.Item(3, 8).Value()  = -1.350072057646098
$ws2.Cells.Item(3, 9).Value()  = "2025-10-14"
$ws2.Cells.Item(3, 10).Value() = "2025-10-14"

# --- Row 4: 标普500 ---
$ws2.Cells.Item(4, 5).Value()  = 2.02
$ws2.Cells.Item(4, 6).Value()  = 3952.6552
$ws2.Cells.Item(4, 7).Value()  = 261.1312556586272
$ws2.Cells.Item(4, 8).Value()  = 7.073806362787042
$ws2.Cells.Item(4, 9).Value()  = "2025-10-14"
$ws2.Cells.Item(4, 10).Value() = "2025-10-14"

# --- Row 5: 纳斯达克 ---
$ws2.Cells.Item(5, 5).Value()  = 5.61
$ws2.Cells.Item(5, 6).Value()  = 2938.1814
$ws2.Cells.Item(5, 7).Value()  = 189.0746113706732
$ws2.Cells.Item(5, 8).Value()  = 6.877674310532826
$ws2.Cells.Item(5, 9).Value()  = "2025-10-14"
$ws2.Cells.Item(5, 10).Value() = "2025-10-14"

# --- Row 6: 兴全合宜混合A ---
$ws2.Cells.Item(6, 5).Value()  = 1.9104
$ws2.Cells.Item(6, 6).Value()  = 7378.174944
$ws2.Cells.Item(6, 7).Value()  = 177.9410835348826
$ws2.Cells.Item(6, 8).Value()  = 2.471323667859145
$ws2.Cells.Item(6, 9).Value()  = "2025-10-14"
$ws2.Cells.Item(6, 10).Value() = "2025-10-14"

# --- New holding rows ---
function Set-HoldRow($row, $asset, $qty, $avgCost, $totalCost, $price, $mktVal, $pnl, $pnlPct, $priceDate, $updDate) {
    $ws2.Cells.Item($row, 1).Value()  = $asset
    $ws2.Cells.Item($row, 2).Value()  = $qty
    $ws2.Cells.Item($row, 3).Value()  = $avgCost
    $ws2.Cells.Item($row, 4).Value()  = $totalCost
    $ws2.Cells.Item($row, 5).Value()  = $price
    $ws2.Cells.Item($row, 6).Value()  = $mktVal
    $ws2.Cells.Item($row, 7).Value()  = $pnl
    $ws2.Cells.Item($row, 8).Value()  = $pnlPct
    $ws2.Cells.Item($row, 9).Value()  = $priceDate
    $ws2.Cells.Item($row, 10).Value() = $updDate
}

Set-HoldRow 7  "有色etf" 1800 0.9037777777777778 1626.8 0.898 1616.4 -10.39999999999986    -0.639291861322834   "2025-10-14" "2025-10-14"
Set-HoldRow 8  "芯片etf" 5500 1.089545454545455  5992.5 1.04  5720   -272.5                -4.547350855235711   "2025-10-14" "2025-10-14"
Set-HoldRow 9  "科创50"   2500 1.522              3805   1.481 3702.5 -102.4999999999995    -2.693823915900119   "2025-10-14" "2025-10-14"
Set-HoldRow 10 "稀土etf" 800  1.36125             1089   1.361 1088.8 -0.2000000000000455  -0.01836547291093163 "2025-10-14" "2025-10-14"

Write-Output "edit complete"
